$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 26.81310566666667
$ws.Range("H2").Value2 = 80.439317
$ws.Range("I2").Value2 = 0.004518206005002021
$ws.Range("J2").Value2 = 0.004518206005002021
$ws.Range("M2").Value2 = 576.300578
$ws.Range("N2").Value2 = 1728.901734
$ws.Range("O2").Value2 = 0.8614732012478776
$ws.Range("P2").Value2 = 0.8614732012478775
$ws.Range("Q2").Value2 = 15452.40829367508
$ws.Range("R2").Value2 = 139071.6746430757
$ws.Range("S2").Value2 = 0.003892313391026475
$ws.Range("T2").Value2 = 0.003892313391026475
$ws.Range("G3").Value2 = 26.81310566666667
$ws.Range("H3").Value2 = 80.439317
$ws.Range("I3").Value2 = 0.004518206005002021
$ws.Range("J3").Value2 = 0.004518206005002021
$ws.Range("O3").Value2 = 0.001786049553652741
$ws.Range("P3").Value2 = 0.001786049553652741
$ws.Range("Q3").Value2 = 32.03670978482022
$ws.Range("R3").Value2 = 288.330388063382
$ws.Range("S3").Value2 = 0.000008069739818544993
$ws.Range("T3").Value2 = 0.000008069739818544991
$ws.Range("G4").Value2 = 26.81310566666667
$ws.Range("H4").Value2 = 80.439317
$ws.Range("I4").Value2 = 0.004518206005002021
$ws.Range("J4").Value2 = 0.004518206005002021
$ws.Range("M4").Value2 = 34.99993866666667
$ws.Range("N4").Value2 = 104.999816
$ws.Range("O4").Value2 = 0.05231906813505349
$ws.Range("P4").Value2 = 0.05231906813505348
$ws.Range("Q4").Value2 = 938.4570537961858
$ws.Range("R4").Value2 = 8446.113484165671
$ws.Range("S4").Value2 = 0.0002363883278239085
$ws.Range("T4").Value2 = 0.0002363883278239085
$ws.Range("G5").Value2 = 26.81310566666667
$ws.Range("H5").Value2 = 80.439317
$ws.Range("I5").Value2 = 0.004518206005002021
$ws.Range("J5").Value2 = 0.004518206005002021
$ws.Range("M5").Value2 = 56.47565533333333
$ws.Range("N5").Value2 = 169.426966
$ws.Range("O5").Value2 = 0.08442168106341624
$ws.Range("P5").Value2 = 0.08442168106341623
$ws.Range("Q5").Value2 = 1514.287714046914
$ws.Range("R5").Value2 = 13628.58942642222
$ws.Range("S5").Value2 = 0.0003814345463330926
$ws.Range("T5").Value2 = 0.0003814345463330926
$ws.Range("G6").Value2 = 5771.873535333333
$ws.Range("I6").Value2 = 0.9726032482643521
$ws.Range("J6").Value2 = 0.9726032482643523
$ws.Range("M6").Value2 = 576.300578
$ws.Range("N6").Value2 = 1728.901734
$ws.Range("O6").Value2 = 0.8614732012478776
$ws.Range("P6").Value2 = 0.8614732012478775
$ws.Range("Q6").Value2 = 3326334.054555503
$ws.Range("R6").Value2 = 29937006.49099953
$ws.Range("S6").Value2 = 0.8378716338263756
$ws.Range("T6").Value2 = 0.8378716338263758
$ws.Range("G7").Value2 = 5771.873535333333
$ws.Range("I7").Value2 = 0.9726032482643521
$ws.Range("J7").Value2 = 0.9726032482643523
$ws.Range("O7").Value2 = 0.001786049553652741
$ws.Range("P7").Value2 = 0.001786049553652741
$ws.Range("Q7").Value2 = 6896.323002077142
$ws.Range("R7").Value2 = 62066.90701869428
$ws.Range("S7").Value2 = 0.001737117597443752
$ws.Range("T7").Value2 = 0.001737117597443752
$ws.Range("G8").Value2 = 5771.873535333333
$ws.Range("I8").Value2 = 0.9726032482643521
$ws.Range("J8").Value2 = 0.9726032482643523
$ws.Range("M8").Value2 = 34.99993866666667
$ws.Range("N8").Value2 = 104.999816
$ws.Range("O8").Value2 = 0.05231906813505349
$ws.Range("P8").Value2 = 0.05231906813505348
$ws.Range("Q8").Value2 = 202015.2197284232
$ws.Range("R8").Value2 = 1818136.977555808
$ws.Range("S8").Value2 = 0.05088569561431697
$ws.Range("T8").Value2 = 0.05088569561431698
$ws.Range("G9").Value2 = 5771.873535333333
$ws.Range("I9").Value2 = 0.9726032482643521
$ws.Range("J9").Value2 = 0.9726032482643523
$ws.Range("M9").Value2 = 56.47565533333333
$ws.Range("N9").Value2 = 169.426966
$ws.Range("O9").Value2 = 0.08442168106341624
$ws.Range("P9").Value2 = 0.08442168106341623
$ws.Range("Q9").Value2 = 325970.3404090735
$ws.Range("R9").Value2 = 2933733.063681661
$ws.Range("S9").Value2 = 0.08210880122621578
$ws.Range("T9").Value2 = 0.08210880122621578
$ws.Range("G10").Value2 = 132.4457753333333
$ws.Range("H10").Value2 = 397.337326
$ws.Range("I10").Value2 = 0.02231808970163987
$ws.Range("J10").Value2 = 0.02231808970163988
$ws.Range("M10").Value2 = 576.300578
$ws.Range("N10").Value2 = 1728.901734
$ws.Range("O10").Value2 = 0.8614732012478776
$ws.Range("P10").Value2 = 0.8614732012478775
$ws.Range("Q10").Value2 = 76328.57687825814
$ws.Range("R10").Value2 = 686957.1919043233
$ws.Range("S10").Value2 = 0.01922643618100899
$ws.Range("T10").Value2 = 0.01922643618100899
$ws.Range("G11").Value2 = 132.4457753333333
$ws.Range("H11").Value2 = 397.337326
$ws.Range("I11").Value2 = 0.02231808970163987
$ws.Range("J11").Value2 = 0.02231808970163988
$ws.Range("O11").Value2 = 0.001786049553652741
$ws.Range("P11").Value2 = 0.001786049553652741
$ws.Range("Q11").Value2 = 158.2482432034884
$ws.Range("R11").Value2 = 1424.234188831396
$ws.Range("S11").Value2 = 0.00003986121414999573
$ws.Range("T11").Value2 = 0.00003986121414999573
$ws.Range("G12").Value2 = 132.4457753333333
$ws.Range("H12").Value2 = 397.337326
$ws.Range("I12").Value2 = 0.02231808970163987
$ws.Range("J12").Value2 = 0.02231808970163988
$ws.Range("M12").Value2 = 34.99993866666667
$ws.Range("N12").Value2 = 104.999816
$ws.Range("O12").Value2 = 0.05231906813505349
$ws.Range("P12").Value2 = 0.05231906813505348
$ws.Range("Q12").Value2 = 4635.594013325779
$ws.Range("R12").Value2 = 41720.34611993202
$ws.Range("S12").Value2 = 0.001167661655744332
$ws.Range("T12").Value2 = 0.001167661655744332
$ws.Range("G13").Value2 = 132.4457753333333
$ws.Range("H13").Value2 = 397.337326
$ws.Range("I13").Value2 = 0.02231808970163987
$ws.Range("J13").Value2 = 0.02231808970163988
$ws.Range("M13").Value2 = 56.47565533333333
$ws.Range("N13").Value2 = 169.426966
$ws.Range("O13").Value2 = 0.08442168106341624
$ws.Range("P13").Value2 = 0.08442168106341623
$ws.Range("Q13").Value2 = 7479.961958081434
$ws.Range("R13").Value2 = 67319.65762273292
$ws.Range("S13").Value2 = 0.001884130650736556
$ws.Range("T13").Value2 = 0.001884130650736556
$ws.Range("G14").Value2 = 3.326003
$ws.Range("H14").Value2 = 9.978009
$ws.Range("I14").Value2 = 0.0005604560290058679
$ws.Range("J14").Value2 = 0.000560456029005868
$ws.Range("M14").Value2 = 576.300578
$ws.Range("N14").Value2 = 1728.901734
$ws.Range("O14").Value2 = 0.8614732012478776
$ws.Range("P14").Value2 = 0.8614732012478775
$ws.Range("Q14").Value2 = 1916.777451329734
$ws.Range("R14").Value2 = 17250.99706196761
$ws.Range("S14").Value2 = 0.0004828178494663584
$ws.Range("T14").Value2 = 0.0004828178494663584
$ws.Range("G15").Value2 = 3.326003
$ws.Range("H15").Value2 = 9.978009
$ws.Range("I15").Value2 = 0.0005604560290058679
$ws.Range("J15").Value2 = 0.000560456029005868
$ws.Range("O15").Value2 = 0.001786049553652741
$ws.Range("P15").Value2 = 0.001786049553652741
$ws.Range("Q15").Value2 = 3.973959383112667
$ws.Range("R15").Value2 = 35.765634448014
$ws.Range("S15").Value2 = 0.000001001002240447918
$ws.Range("T15").Value2 = 0.000001001002240447918
$ws.Range("G16").Value2 = 3.326003
$ws.Range("H16").Value2 = 9.978009
$ws.Range("I16").Value2 = 0.0005604560290058679
$ws.Range("J16").Value2 = 0.000560456029005868
$ws.Range("M16").Value2 = 34.99993866666667
$ws.Range("N16").Value2 = 104.999816
$ws.Range("O16").Value2 = 0.05231906813505349
$ws.Range("P16").Value2 = 0.05231906813505348
$ws.Range("Q16").Value2 = 116.4099010051493
$ws.Range("R16").Value2 = 1047.689109046344
$ws.Range("S16").Value2 = 0.00002932253716825951
$ws.Range("T16").Value2 = 0.00002932253716825952
$ws.Range("G17").Value2 = 3.326003
$ws.Range("H17").Value2 = 9.978009
$ws.Range("I17").Value2 = 0.0005604560290058679
$ws.Range("J17").Value2 = 0.000560456029005868
$ws.Range("M17").Value2 = 56.47565533333333
$ws.Range("N17").Value2 = 169.426966
$ws.Range("O17").Value2 = 0.08442168106341624
$ws.Range("P17").Value2 = 0.08442168106341623
$ws.Range("Q17").Value2 = 187.8381990656327
$ws.Range("R17").Value2 = 1690.543791590694
$ws.Range("S17").Value2 = 0.00004731464013080214
$ws.Range("T17").Value2 = 0.00004731464013080214
